{"js": "// Swap the two halves of the \"Desempenho da Classifica\u00e7\u00e3o SITS \u2013 Resultados\"\n// bullet so it reads \"Resultados \u2013 Desempenho da Classifica\u00e7\u00e3o SITS\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst oldText = \"Desempenho da Classifica\u00e7\u00e3o SITS \u2013 Resultados\";\nconst newText = \"Resultados \u2013 Desempenho da Classifica\u00e7\u00e3o SITS\";\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text === oldText) {\n    target = para;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the target paragraph: \" + oldText);\n}\n\n// Replace the paragraph's text in place, preserving the run formatting\n// (all runs in this paragraph already share the same Times New Roman rPr).\ntarget.insertText(newText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Swap the two halves of the \"Desempenho da Classifica\u00e7\u00e3o SITS \u2013 Resultados\"\n# bullet so it reads \"Resultados \u2013 Desempenho da Classifica\u00e7\u00e3o SITS\".\n$d = $word.ActiveDocument\n\n$oldText = \"Desempenho da Classifica\u00e7\u00e3o SITS \u2013 Resultados\"\n$newText = \"Resultados \u2013 Desempenho da Classifica\u00e7\u00e3o SITS\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 1)\n\nif (-not $found) {\n    throw \"Could not find the target text: $oldText\"\n}\n"}
